$wb = $excel.ActiveWorkbook

# Rename the existing sheet "Munka1" to "sütemények"
$ws1 = $wb.Worksheets.Item("Munka1")
$ws1.Name = "sütemények"

# Add new sheet "kávék" after the existing sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "kávék"

# Fill header row
$ws2.Range("A1").Value = "név"
$ws2.Range("B1").Value = "kép"
$ws2.Range("C1").Value = "kávé"
$ws2.Range("D1").Value = "víz/tej"
$ws2.Range("E1").Value = "tejszín"
$ws2.Range("F1").Value = "alkotás"
$ws2.Range("G1").Value = "mekkora?"

# Add a comment to G1
$ws2.Range("G1").AddComment("Maga a kávé mekkora?")

# Activate the new sheet (kávék) as the selected / visible tab
$ws2.Activate()
$ws2.Range("L7").Select()
